$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H26").Value = 15600
$ws.Range("J26").Value = 15600
$ws.Range("L26").Value = 15600
$ws.Range("N26").Value = -16288

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 3088.0557
$ws.Range("I62").Value = 2100.5293
$ws.Range("J62").Value = 19876
$ws.Range("K62").Value = 2100.5293
$ws.Range("L62").Value = 19876
$ws.Range("M62").Value = -1476.5293
$ws.Range("N62").Value = -21124

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H65").Value = 3088.0557
$ws.Range("I65").Value = 2100.5293
$ws.Range("J65").Value = 19876
$ws.Range("K65").Value = 10502.6465
$ws.Range("L65").Value = 99380
$ws.Range("M65").Value = -7382.646500000001
$ws.Range("N65").Value = -105620

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 3380.842
$ws.Range("I138").Value = 1777.7646
$ws.Range("J138").Value = 3842.7458
$ws.Range("K138").Value = 5333.293799999999
$ws.Range("L138").Value = 11528.2374
$ws.Range("M138").Value = -193.2937999999995
$ws.Range("N138").Value = -21808.2374

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 21242.018
$ws.Range("I32").Value = 22639.633
$ws.Range("J32").Value = 7545.4
$ws.Range("K32").Value = 22639.633
$ws.Range("L32").Value = 7545.4
$ws.Range("M32").Value = -22352.633
$ws.Range("N32").Value = -8119.4

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2927.4443
$ws.Range("I132").Value = 2610.5217
$ws.Range("J132").Value = 4749.75
$ws.Range("K132").Value = 7831.5651
$ws.Range("L132").Value = 14249.25
$ws.Range("M132").Value = -5301.5651
$ws.Range("N132").Value = -19309.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1272.9375
$ws.Range("I94").Value = 1185.5385
$ws.Range("K94").Value = 1185.5385
$ws.Range("M94").Value = -734.5385000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2640.842
$ws.Range("I99").Value = 869.7143
$ws.Range("J99").Value = 7600
$ws.Range("K99").Value = 869.7143
$ws.Range("L99").Value = 7600
$ws.Range("M99").Value = 628.2857
$ws.Range("N99").Value = -10596

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 4045.7334
$ws.Range("I105").Value = 3587.6667
$ws.Range("J105").Value = 4732.8335
$ws.Range("K105").Value = 3587.6667
$ws.Range("L105").Value = 4732.8335
$ws.Range("M105").Value = -1840.6667
$ws.Range("N105").Value = -8226.833500000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2166146.5
$ws.Range("I58").Value = 4330437.5
$ws.Range("J58").Value = 1855.7142
$ws.Range("K58").Value = 4330437.5
$ws.Range("L58").Value = 1855.7142
$ws.Range("M58").Value = -4330234.5
$ws.Range("N58").Value = -2261.7142

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H59").Value = 25607.055
$ws.Range("J59").Value = 25607.055
$ws.Range("L59").Value = 25607.055
$ws.Range("N59").Value = -27897.055

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2270.5356
$ws.Range("I132").Value = 2072.8838
$ws.Range("J132").Value = 2924.3076
$ws.Range("K132").Value = 6218.651400000001
$ws.Range("L132").Value = 8772.9228
$ws.Range("M132").Value = -3688.651400000001
$ws.Range("N132").Value = -13832.9228

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 2166146.5
$ws.Range("I136").Value = 4330437.5
$ws.Range("J136").Value = 1855.7142
$ws.Range("K136").Value = 12991312.5
$ws.Range("L136").Value = 5567.142599999999
$ws.Range("M136").Value = -12988762.5
$ws.Range("N136").Value = -10667.1426

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 132.85715
$ws.Range("I7").Value = 87.25
$ws.Range("K7").Value = 261.75
$ws.Range("M7").Value = -149.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H25").Value = 4000
$ws.Range("I25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("M25").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H30").Value = 4000
$ws.Range("I30").Value = 0
$ws.Range("K30").Value = 0
$ws.Range("M30").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H69").Value = 166669360
$ws.Range("J69").Value = 250003500
$ws.Range("L69").Value = 750010500
$ws.Range("N69").Value = -750012122

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H72").Value = 166669360
$ws.Range("J72").Value = 250003500
$ws.Range("L72").Value = 2250031500
$ws.Range("N72").Value = -2250039612

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H109").Value = 2256.1667
$ws.Range("I109").Value = 633
$ws.Range("J109").Value = 3289.0908
$ws.Range("K109").Value = 1899
$ws.Range("L109").Value = 9867.2724
$ws.Range("M109").Value = -859
$ws.Range("N109").Value = -11947.2724

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H112").Value = 2363.5
$ws.Range("I112").Value = 1484.6666
$ws.Range("J112").Value = 5000
$ws.Range("K112").Value = 4453.9998
$ws.Range("L112").Value = 15000
$ws.Range("M112").Value = -3345.9998
$ws.Range("N112").Value = -17216

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H118").Value = 13509.857
$ws.Range("I118").Value = 2428.1667
$ws.Range("J118").Value = 80000
$ws.Range("K118").Value = 7284.500100000001
$ws.Range("L118").Value = 240000
$ws.Range("M118").Value = -6041.500100000001
$ws.Range("N118").Value = -242486

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1257.1464
$ws.Range("J131").Value = 1073.3226
$ws.Range("L131").Value = 3219.9678
$ws.Range("N131").Value = -13299.9678

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H14").Value = 175167020
$ws.Range("I14").Value = 175167020
$ws.Range("K14").Value = 175167020
$ws.Range("M14").Value = -175166852

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H49").Value = 19800
$ws.Range("J49").Value = 19800
$ws.Range("L49").Value = 19800
$ws.Range("N49").Value = -20168

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5397.7334
$ws.Range("I70").Value = 4988.815
$ws.Range("K70").Value = 4988.815
$ws.Range("M70").Value = -4718.815

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 5397.7334
$ws.Range("I73").Value = 4988.815
$ws.Range("K73").Value = 4988.815
$ws.Range("M73").Value = -4052.815

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2743.8
$ws.Range("I126").Value = 1982.0769
$ws.Range("J126").Value = 3569
$ws.Range("K126").Value = 5946.2307
$ws.Range("L126").Value = 10707
$ws.Range("M126").Value = -3476.2307
$ws.Range("N126").Value = -15647

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 9097
$ws.Range("I16").Value = 916.4
$ws.Range("J16").Value = 50000
$ws.Range("K16").Value = 916.4
$ws.Range("L16").Value = 50000
$ws.Range("M16").Value = -746.4
$ws.Range("N16").Value = -50340

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 11319.429
$ws.Range("I132").Value = 14611.556
$ws.Range("J132").Value = 5393.6
$ws.Range("K132").Value = 43834.66800000001
$ws.Range("L132").Value = 16180.8
$ws.Range("M132").Value = -41304.66800000001
$ws.Range("N132").Value = -21240.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H24").Value = 29999.666
$ws.Range("J24").Value = 29999.666
$ws.Range("L24").Value = 29999.666
$ws.Range("N24").Value = -30459.666

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2184.3872
$ws.Range("I132").Value = 2142.2083
$ws.Range("J132").Value = 2329
$ws.Range("K132").Value = 6426.624899999999
$ws.Range("L132").Value = 6987
$ws.Range("M132").Value = -3896.624899999999
$ws.Range("N132").Value = -12047
